$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.544.04'
$ws.Range("E2").Value = '  +5.40%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.724.20'
$ws.Range("E3").Value = '  +4.22%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '226.13'
$ws.Range("E5").Value = '  +3.48%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5378'
$ws.Range("E6").Value = '  +2.69%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.003'
$ws.Range("E7").Value = '  -0.06%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2677'
$ws.Range("E8").Value = '  +0.73%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06616'
$ws.Range("E9").Value = '  +4.18%  '

$ws.Range("E10").Value = '  +6.31%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07728'
$ws.Range("E11").Value = '  +0.55%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.618'
$ws.Range("E12").Value = '  -0.11%  '

$ws.Range("B13").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C13").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.962.48'
$ws.Range("E13").Value = '  +4.24%  '

$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.702.55'
$ws.Range("E14").Value = '  +4.26%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5883'
$ws.Range("E15").Value = '  +4.92%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0₅8319'
$ws.Range("E16").Value = '  +1.69%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '68.04'
$ws.Range("E17").Value = '  +4.02%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '27.561.83'
$ws.Range("E18").Value = '  +5.47%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '222.31'
$ws.Range("E19").Value = '  +15.71%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.003'
$ws.Range("E20").Value = '  +0.02%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.744'
$ws.Range("E21").Value = '  +2.03%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.72'
$ws.Range("E22").Value = '  +1.82%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.109'
$ws.Range("E23").Value = '  +2.65%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.003'
$ws.Range("E24").Value = '  -0.04%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '148.26'
$ws.Range("E25").Value = '  +2.43%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.696'
$ws.Range("E26").Value = '  +12.20%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1236'
$ws.Range("E27").Value = '  +3.46%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.408'
$ws.Range("E28").Value = '  +2.01%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '16.70'
$ws.Range("E29").Value = '  +4.76%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05541'
$ws.Range("E30").Value = '  +1.86%  '

$ws.Range("E31").Value = '  +2.49%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.560'
$ws.Range("E32").Value = '  +2.82%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.466'
$ws.Range("E33").Value = '  +3.00%  '

$ws.Range("E34").Value = '  +6.53%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9640'
$ws.Range("E35").Value = '  +1.39%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.821'
$ws.Range("E36").Value = '  +1.52%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.443'
$ws.Range("E37").Value = '  +1.71%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5966'
$ws.Range("E38").Value = '  +5.32%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01648'
$ws.Range("E39").Value = '  +4.16%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.942'
$ws.Range("E40").Value = '  +1.26%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.060.17'
$ws.Range("E41").Value = '  +3.16%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8535'
$ws.Range("E42").Value = '  +2.48%  '

$ws.Range("E43").Value = '  -0.02%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '101.58'
$ws.Range("E44").Value = '  +0.16%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.867.61'
$ws.Range("E45").Value = '  +4.13%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0₈114'
$ws.Range("E46").Value = '  +17.94%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '59.14'
$ws.Range("E47").Value = '  +2.49%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.230'
$ws.Range("E48").Value = '  +2.73%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.4438'
$ws.Range("E49").Value = '  +2.26%  '

$ws.Range("E50").Value = '  +0.27%  '

$ws.Range("E51").Value = '  +1.61%  '
